$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the Age (years) minimum value: 18 -> 17
$ws.Range("D2").Value = 17

# Insert a new row for "Serum glucose" right before the "Serum bicarbonate" row (row 17)
$ws.Rows.Item(17).Insert()

# Populate the new row 17 with the Serum glucose data
$ws.Range("A17").Value = "Serum glucose"
$ws.Range("B17").Value = "glucose"
$ws.Range("C17").Value = "mmol/l"
$ws.Range("D17").Value = 0.5
$ws.Range("D17").NumberFormat = "0.0"
$ws.Range("E17").Value = 87

# Clean up the creatinine units label (drop trailing space)
$ws.Range("B19").Value = "creatinine"
